# Updated cryptos list on Thu Mar  2 05:00:06 UTC 2023 with GitHub Actions
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row (rows 2-51) on Sheet1 to the latest scraped values.
#
# Price cells that look like plain decimal numbers (e.g. "1.000", "22.17")
# must be forced to Text so Excel doesn't auto-convert/round them; this is
# done by temporarily setting NumberFormat to "@" before assigning the
# value and then restoring the cell style to "Normal" afterwards so no
# extra formatting is left behind. Prices that already contain more than
# one "." (e.g. "23.540.39") are naturally treated as text by Excel, so
# no special handling is required for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.540.39"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.651.48"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.64%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "300.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3795"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.81"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.16%  "
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.229"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08117"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.435"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.443"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001207"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.84%  "
$ws.Range("D17").Value = "1.659.33"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06983"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.795"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D24").Value = "23.552.24"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.490"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.934"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.238"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "1.843.71"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.016"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.141"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.034"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02749"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08717"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.021"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2456"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "13.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06901"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6945"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.326"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6460"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.277"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.932"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07822"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.178"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.64%  "
